# Refresh the cryptos list: update "Price" (D) and "Volume(1h)" (E) columns.
# Price cells are forced to Text (leading apostrophe) so strings such as
# "7.50" / "1.00" / "14.40" keep their trailing zeros instead of being
# auto-coerced to numbers 7.5 / 1 / 14.4 by the COM value-setter; the style
# is then reset to "Normal" so no stray quote-prefix formatting is left on
# the cell (matching the source file, where these are plain inline strings).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.Value = "'63.724.82"
$c.Style = 'Normal'
$ws.Range('E2').Value = '  -2.39%  '
$c = $ws.Range('D3')
$c.Value = "'3.049.99"
$c.Style = 'Normal'
$ws.Range('E3').Value = '  -2.06%  '
$ws.Range('E4').Value = '  +0.00%  '
$c = $ws.Range('D5')
$c.Value = "'555.88"
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -1.00%  '
$c = $ws.Range('D6')
$c.Value = "'141.52"
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -2.81%  '
$ws.Range('E7').Value = '  +0.04%  '
$c = $ws.Range('D8')
$c.Value = "'3.048.48"
$c.Style = 'Normal'
$ws.Range('E8').Value = '  -1.95%  '
$ws.Range('E9').Value = '  +3.54%  '
$c = $ws.Range('D10')
$c.Value = "'0.152"
$c.Style = 'Normal'
$ws.Range('E10').Value = '  -0.79%  '
$c = $ws.Range('D11')
$c.Value = "'6.18"
$c.Style = 'Normal'
$ws.Range('E11').Value = '  -14.01%  '
$c = $ws.Range('D12')
$c.Value = "'0.477"
$c.Style = 'Normal'
$ws.Range('E12').Value = '  +1.51%  '
$ws.Range('E13').Value = '  -1.47%  '
$c = $ws.Range('D14')
$c.Value = "'35.06"
$c.Style = 'Normal'
$ws.Range('E14').Value = '  -1.64%  '
$c = $ws.Range('D15')
$c.Value = "'3.548.05"
$c.Style = 'Normal'
$ws.Range('E15').Value = '  -1.64%  '
$c = $ws.Range('D16')
$c.Value = "'63.777.87"
$c.Style = 'Normal'
$ws.Range('E16').Value = '  -2.29%  '
$c = $ws.Range('D17')
$c.Value = "'3.045.55"
$c.Style = 'Normal'
$ws.Range('E17').Value = '  -1.84%  '
$ws.Range('E18').Value = '  +0.21%  '
$c = $ws.Range('D19')
$c.Value = "'6.74"
$c.Style = 'Normal'
$ws.Range('E19').Value = '  -1.98%  '
$c = $ws.Range('D20')
$c.Value = "'487.75"
$c.Style = 'Normal'
$ws.Range('E20').Value = '  +0.78%  '
$c = $ws.Range('D21')
$c.Value = "'14.14"
$c.Style = 'Normal'
$ws.Range('E21').Value = '  +1.79%  '
$c = $ws.Range('D22')
$c.Value = "'0.681"
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -0.31%  '
$c = $ws.Range('D23')
$c.Value = "'14.40"
$c.Style = 'Normal'
$ws.Range('E23').Value = '  +6.13%  '
$c = $ws.Range('D24')
$c.Value = "'7.50"
$c.Style = 'Normal'
$c = $ws.Range('D25')
$c.Value = "'82.41"
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +1.18%  '
$ws.Range('E26').Value = '  -0.12%  '
$c = $ws.Range('D27')
$c.Value = "'2.78"
$c.Style = 'Normal'
$ws.Range('E27').Value = '  -0.57%  '
$c = $ws.Range('D28')
$c.Value = "'8.07"
$c.Style = 'Normal'
$ws.Range('E28').Value = '  -2.17%  '
$ws.Range('E29').Value = '  -2.04%  '
$c = $ws.Range('D30')
$c.Value = "'1.00"
$c.Style = 'Normal'
$ws.Range('E30').Value = '  +0.24%  '
$c = $ws.Range('D31')
$c.Value = "'26.23"
$c.Style = 'Normal'
$ws.Range('E31').Value = '  -0.09%  '
$ws.Range('E32').Value = '  -0.75%  '
$c = $ws.Range('D33')
$c.Value = "'2.44"
$c.Style = 'Normal'
$ws.Range('E33').Value = '  -2.30%  '
$c = $ws.Range('D34')
$c.Value = "'5.65"
$c.Style = 'Normal'
$ws.Range('E34').Value = '  -1.05%  '
$c = $ws.Range('D35')
$c.Value = "'6.18"
$c.Style = 'Normal'
$ws.Range('E35').Value = '  -1.22%  '
$c = $ws.Range('D36')
$c.Value = "'55.16"
$c.Style = 'Normal'
$ws.Range('E36').Value = '  -0.17%  '
$c = $ws.Range('D37')
$c.Value = "'0.0407"
$c.Style = 'Normal'
$ws.Range('E37').Value = '  -1.46%  '
$c = $ws.Range('D38')
$c.Value = "'440.08"
$c.Style = 'Normal'
$ws.Range('E38').Value = '  -6.83%  '
$c = $ws.Range('D39')
$c.Value = "'0.0812"
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -2.71%  '
$c = $ws.Range('D40')
$c.Value = "'3.017.47"
$c.Style = 'Normal'
$ws.Range('E40').Value = '  -0.23%  '
$c = $ws.Range('D41')
$c.Value = "'2.75"
$c.Style = 'Normal'
$ws.Range('E41').Value = '  -6.35%  '
$c = $ws.Range('D42')
$c.Value = "'8.30"
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -0.11%  '
$ws.Range('E43').Value = '  -1.19%  '
$c = $ws.Range('D44')
$c.Value = "'0.269"
$c.Style = 'Normal'
$ws.Range('E44').Value = '  +3.53%  '
$c = $ws.Range('D45')
$c.Value = "'27.55"
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -2.99%  '
$c = $ws.Range('D46')
$c.Value = "'2.20"
$c.Style = 'Normal'
$ws.Range('E46').Value = '  +2.54%  '
$ws.Range('E48').Value = '  +0.19%  '
$c = $ws.Range('D49')
$c.Value = "'117.69"
$c.Style = 'Normal'
$ws.Range('E49').Value = '  +0.65%  '
$c = $ws.Range('D50')
$c.Value = "'0.0₃0511"
$c.Style = 'Normal'
$ws.Range('E50').Value = '  -2.96%  '
$ws.Range('E51').Value = '  -0.91%  '
